$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
  "A1","B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","O1","P1","Q1","R1","S1","T1","U1","V1","W1","X1",
  "A2","B2","C2","D2","E2","F2","G2","H2","I2","J2","K2","L2","O2","P2","Q2","R2","S2","T2","U2","V2","W2","X2",
  "A3","B3","W3","X3",
  "A4","B4","F4","G4","H4","I4","J4","O4","P4","Q4","R4","S4","W4","X4",
  "A5","B5","F5","G5","H5","I5","J5","O5","P5","Q5","R5","S5","W5","X5",
  "A6","B6","F6","G6","H6","I6","J6","O6","P6","Q6","R6","S6","W6","X6",
  "A7","B7","E7","F7","G7","H7","I7","J7","N7","O7","P7","Q7","R7","S7","W7","X7",
  "A8","B8","E8","F8","G8","H8","I8","J8","N8","O8","P8","Q8","R8","S8","W8","X8",
  "A9","B9","W9","X9",
  "A10","B10","N10","O10","P10","Q10","R10","S10","T10","W10","X10",
  "A11","B11","N11","O11","P11","Q11","R11","S11","T11","W11","X11",
  "A12","B12","F12","G12","H12","I12","J12","N12","O12","P12","Q12","R12","S12","T12","W12","X12",
  "A13","B13","N13","O13","P13","Q13","R13","S13","T13","W13","X13",
  "A14","B14","N14","O14","P14","Q14","R14","S14","T14","W14","X14",
  "A15","B15","F15","W15","X15",
  "A16","B16","W16","X16",
  "A17","B17","N17","O17","P17","Q17","R17","S17","W17","X17",
  "A18","B18","W18","X18",
  "A19","B19","W19","X19",
  "A20","B20","C20","D20","E20","F20","S20","T20","U20","V20","W20","X20"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = "x"
}
